$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the distribution-type code in column F (rows 34-39) from "PT,(E-1B)" to "RT,(E-1B)"
$ws.Range("F34:F39").Value = "RT,(E-1B)"
